$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-7 from 45174 to 45175
$ws.Range("C2:C7").Value = 45175
